# Weekly update of "Fruta, Terminal Hortofrutícola Agro Chillán - Pera" data.
# A new week's pair of rows (Especial / Primera) is inserted at the top of the
# date-ordered block (rows 129-130), pushing the existing history down by two
# rows. The two rows that previously occupied 129/130 are preserved by being
# copied into the newly opened 131/132 slots (everything below cascades by
# the same +2 offset), and the data that falls off the bottom of the original
# A1:T155 range reappears as two brand-new rows 156/157.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Open up two fresh rows right after the current top data rows (129, 130).
#    Everything from old row 131 downward shifts to row+2.
$ws.Rows("131:132").Insert()

# 2) The rows that used to be at 129/130 now need to live on (shifted by two)
#    immediately below the new entry, i.e. at 131/132.
$ws.Range("A129:T129").Copy($ws.Range("A131:T131"))
$ws.Range("A130:T130").Copy($ws.Range("A132:T132"))

# 3) Overwrite row 129/130 with this week's new readings.
#    Row 129 - "Especial" quality.
$ws.Range("D129").Value = 44505
$ws.Range("M129").Value = 60
$ws.Range("N129").Value = 10500
$ws.Range("O129").Value = 11000
$ws.Range("P129").Value = 10750
$ws.Range("S129").Value = 672

# Row 130 - "Primera" quality: only the date moves to the new reporting day.
$ws.Range("D130").Value = 44505
